$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly data point was added for this market/product. It lands at
# row 128, pushing the existing data rows 128:190 down to 129:191 (the
# sheet's used range grows from A1:R190 to A1:R191).
$ws.Rows(128).Insert()

$ws.Range("A128").Value = 4
$ws.Range("B128").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C128").Value = "Los Lagos"
$ws.Range("D128").Value = 44466
$ws.Range("E128").Value = 10
$ws.Range("F128").Value = 100112023
$ws.Range("G128").Value = "Brócoli"
$ws.Range("H128").Value = "Sin especificar"
$ws.Range("I128").Value = "Segunda"
$ws.Range("J128").Value = 750
$ws.Range("K128").Value = 1000
$ws.Range("L128").Value = 1000
$ws.Range("M128").Value = 1000
$ws.Range("N128").Value = "$/unidad"
$ws.Range("O128").Value = "Región del Maule"
$ws.Range("P128").Value = 1000
$ws.Range("Q128").Value = 1
$ws.Range("R128").Value = "Hortaliza"
